$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Reproduce the formatting of the new column N (rows 2-33) by copying it
#    from column M, which already carries the exact same style pattern one
#    column to the left (header border row, year header, data rows, total
#    row). This reuses the existing style records instead of creating new
#    (duplicate) ones.
$ws.Range("M2:M33").Copy()
$ws.Range("N2:N33").PasteSpecial(-4122)

# 2) Fill in the new year column (2022) of data.
$ws.Range("N3").Value = 2022

$ws.Range("N4").Value = 11.927942610539198
$ws.Range("N5").Value = 3.0909744679837434
$ws.Range("N6").Value = 20.963679772397647

$ws.Range("N7").Value = 4.6002717699014832
$ws.Range("N8").Value = 0
$ws.Range("N9").Value = 9.112830865859129

$ws.Range("N10").Value = 3.5391993253978327
$ws.Range("N11").Value = 0.30955295909412422
$ws.Range("N12").Value = 6.73157537222552

$ws.Range("N13").Value = 3.9173330796393815
$ws.Range("N14").Value = 0.7444796831494469
$ws.Range("N15").Value = 7.104530072727953

$ws.Range("N16").Value = 23.0957399744971
$ws.Range("N17").Value = 2.6274648905004008
$ws.Range("N18").Value = 43.176223433734158

$ws.Range("N19").Value = 7.6660105666632132
$ws.Range("N20").Value = 0.83437630371297455
$ws.Range("N21").Value = 14.406256431364477

$ws.Range("N22").Value = 34.201612992199827
$ws.Range("N23").Value = 4.4521615244201058
$ws.Range("N24").Value = 63.433733622066185

$ws.Range("N25").Value = 20.535408979625672
$ws.Range("N26").Value = 7.8632542639432348
$ws.Range("N27").Value = 33.368028499329796

$ws.Range("N28").Value = 19.301652062045072
$ws.Range("N29").Value = 7.1220113855063829
$ws.Range("N30").Value = 34.008685896558866

$ws.Range("N31").Value = 7.8668258762379715
$ws.Range("N32").Value = 1.7266187050359711
$ws.Range("N33").Value = 13.723068478111704

# 3) Row 34 is the blank spacer row right below the table. Give N34 its own
#    plain "Times New Roman 11" look (matching the rest of the workbook's
#    body font) instead of inheriting the table column's centred alignment.
$n34 = $ws.Range("N34")
$n34.Font.Name = "Times New Roman"
$n34.Font.Size = 11
$n34.VerticalAlignment = -4107

# 4) Update the remembered selection, as recorded in the sheet view.
$ws.Range("O6").Select()
